$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ---
# Header row: Flashcard button label renamed, new "Notes" button header
$ws.Range("F1").Value = "Button_Flashcard"
$ws.Range("G1").Value = "Button_Notes"

# Row 9: lesson name expanded, and the stray flashcard link removed
$ws.Range("C9").Value = "Market Around Us"
$ws.Range("G9").ClearContents()

# Row 8: flashcard link target updated to new html file
$ws.Range("F8").Value = "/FlashCardTest.html"

# --- Column widths ---
# (input values chosen so the saved OOXML "width" lands as close as possible
# to the target 66.88671875 / 39.6640625 / 12.77734375 / 23.5546875)
$ws.Columns("D").ColumnWidth = 65.9167
$ws.Columns("E").ColumnWidth = 38.7501
$ws.Columns("F").ColumnWidth = 11.9167
$ws.Columns("G").ColumnWidth = 22.5834

# --- Selection ---
[void]$ws.Range("F16").Select()
